# Add a new "Slovakia" market sheet to the workbook.
#
# The source data shows the new sheet is a verbatim copy of the existing
# "Germany" sheet (same shared-string references / styles / merges), just
# renamed and placed at the end of the tab strip, which is exactly what
# Excel's "Move or Copy... > Create a copy" does when duplicating a sheet.

$wb = $excel.ActiveWorkbook

# Before duplicating, reproduce Germany's "select all cells" state that
# shows up on the source sheet after the copy operation.
$germany = $wb.Worksheets.Item("Germany")
$germany.Activate()
$germany.Cells.Select()

# Duplicate "Germany" and drop the copy after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$germany.Copy($null, $lastSheet)

# The copy becomes the new last sheet - rename it and give it the
# selection/active-tab state of the new "Slovakia" sheet.
$slovakia = $wb.Worksheets.Item($wb.Worksheets.Count)
$slovakia.Name = "Slovakia"
$slovakia.Range("D18").Select()
